# Fruta / hortaliza, semanal
# Insert a new weekly block (2 rows) at the top of the "Terminal La Palmera
# de La Serena - Caqui" data table (rows 75-76), pushing all existing rows
# (old 75-88) down by two to become rows 77-90. The new rows carry the
# latest week's Mankaki Especial / Primera price observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 75, shifting rows 75:88 (and
# their formatting) down to 77:90.
$ws.Range("A75:A76").EntireRow.Insert()

# New row 75: Mankaki / Especial
$ws.Range("A75").Value = 8
$ws.Range("B75").Value = "Terminal La Palmera de La Serena"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 44694
$ws.Range("E75").Value = 4
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100107
$ws.Range("H75").Value = "Otros"
$ws.Range("I75").Value = 100107001
$ws.Range("J75").Value = "Caqui"
$ws.Range("K75").Value = "Mankaki"
$ws.Range("L75").Value = "Especial"
$ws.Range("M75").Value = 20
$ws.Range("N75").Value = 490000
$ws.Range("O75").Value = 500000
$ws.Range("P75").Value = 495000
$ws.Range("Q75").Value = "$/bins (450 kilos)"
$ws.Range("R75").Value = "Región de O'Higgins"
$ws.Range("S75").Value = 1100
$ws.Range("T75").Value = 450

# New row 76: Mankaki / Primera
$ws.Range("A76").Value = 8
$ws.Range("B76").Value = "Terminal La Palmera de La Serena"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44694
$ws.Range("E76").Value = 4
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100107
$ws.Range("H76").Value = "Otros"
$ws.Range("I76").Value = 100107001
$ws.Range("J76").Value = "Caqui"
$ws.Range("K76").Value = "Mankaki"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 16
$ws.Range("N76").Value = 440000
$ws.Range("O76").Value = 450000
$ws.Range("P76").Value = 445000
$ws.Range("Q76").Value = "$/bins (450 kilos)"
$ws.Range("R76").Value = "Región de O'Higgins"
$ws.Range("S76").Value = 989
$ws.Range("T76").Value = 450
